# Splits three run-on "Programa"/"Bibliografia" paragraphs into separate
# lines joined by manual line breaks (<w:br/>), matching the source diff.
#
# The three paragraphs (PT "Programa", EN "Programa" in italics, and
# "Bibliografia") are each a single run whose <w:t> mashes several
# numbered/sentence items together with no separator. The target splits
# them into one <w:t> per item joined by <w:br/> manual line breaks,
# preserving leading/trailing spaces via xml:space="preserve" exactly as
# in the source.
#
# Plain Find/Replace with a "^l" (manual line break) wildcard token
# reproduces the <w:t>/<w:br/> structure, but this runtime's serializer
# drops the xml:space="preserve" marker on segments that end or begin
# with a space when a <w:br/> immediately follows/precedes - so instead
# each whole paragraph is matched via Find, cleared, and its replacement
# is injected as literal OOXML via Range.InsertXML, which round-trips
# xml:space="preserve" faithfully. Needle/payload text is carried as
# base64 to sidestep any PowerShell quoting issues with accents/quotes.

$d = $word.ActiveDocument

function Replace-ParagraphWithXml($findB64, $xmlB64) {
    $findBytes = [Convert]::FromBase64String($findB64)
    $findText  = [System.Text.Encoding]::UTF8.GetString($findBytes)
    $xmlBytes  = [Convert]::FromBase64String($xmlB64)
    $xmlText   = [System.Text.Encoding]::UTF8.GetString($xmlBytes)

    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for paragraph text starting: $($findText.Substring(0, 40))"
    }
    $rng.Text = "" | Out-Null
    $rng.InsertXML($xmlText) | Out-Null
}

# --- PT Programa ---
Replace-ParagraphWithXml "MSkgSW50cm9kdcOnw6NvIGEgRsOtc2ljYTogc2lzdGVtYXMgZGUgdW5pZGFkZXMsIHJldmlzw6NvIGRlIHZldG9yZXMsIGFuw6FsaXNlIGRpbWVuc2lvbmFsLjIpIENpbmVtw6F0aWNhOiBtb3ZpbWVudG8gdW5pZGltZW5zaW9uYWwsIHF1ZWRhIGxpdnJlLCBtb3ZpbWVudG8gYmlkaW1lbnNpb25hbCwgcHJvasOpdGVpcy4gMykgRGluw6JtaWNhOiBsZWlzIGRlIE5ld3RvbiwgZm9yw6dhcywgZm9yw6dhIGRlIGF0cml0bywgZm9yw6dhIGRlIHJlc2lzdMOqbmNpYSBkbyBhciwgdmVsb2NpZGFkZSB0ZXJtaW5hbCwgbW92aW1lbnRvIGNpcmN1bGFyIHVuaWZvcm1lLCBncmF2aXRhw6fDo28sIGFwbGljYcOnw7Vlcy40KSBFbmVyZ2lhOiB0cmFiYWxobywgZm9yw6dhcyBjb25zZXJ2YXRpdmFzLCBjb25zZXJ2YcOnw6NvIGRlIGVuZXJnaWEgbWVjw6JuaWNhLCBhdHJpdG8sIGFwbGljYcOnw7Vlcy41KSAgTW9tZW50byBsaW5lYXI6IGNlbnRybyBkZSBtYXNzYSwgc2lzdGVtYSBkZSBwYXJ0w61jdWxhcywgY29uc2VydmHDp8OjbyBkbyBtb21lbnRvIGxpbmVhciwgY29saXPDtWVzLCBpbXB1bHNvLjYpIFJvdGHDp8OjbzogdmFyacOhdmVpcyBkbyBtb3ZpbWVudG8gcm90YWNpb25hbCwgZW5lcmdpYSBjaW7DqXRpY2Egcm90YWNpb25hbCwgbW9tZW50byBkZSBpbsOpcmNpYSwgdG9ycXVlLCByb2xhbWVudG8sIGNvbnNlcnZhw6fDo28gZG8gbW9tZW50byBhbmd1bGFyLg==" "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6dD4xKSBJbnRyb2R1w6fDo28gYSBGw61zaWNhOiBzaXN0ZW1hcyBkZSB1bmlkYWRlcywgcmV2aXPDo28gZGUgdmV0b3JlcywgYW7DoWxpc2UgZGltZW5zaW9uYWwuPC93OnQ+PHc6YnIvPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+MikgQ2luZW3DoXRpY2E6IG1vdmltZW50byB1bmlkaW1lbnNpb25hbCwgcXVlZGEgbGl2cmUsIG1vdmltZW50byBiaWRpbWVuc2lvbmFsLCBwcm9qw6l0ZWlzLiA8L3c6dD48dzpici8+PHc6dD4zKSBEaW7Dom1pY2E6IGxlaXMgZGUgTmV3dG9uLCBmb3LDp2FzLCBmb3LDp2EgZGUgYXRyaXRvLCBmb3LDp2EgZGUgcmVzaXN0w6puY2lhIGRvIGFyLCB2ZWxvY2lkYWRlIHRlcm1pbmFsLCBtb3ZpbWVudG8gY2lyY3VsYXIgdW5pZm9ybWUsIGdyYXZpdGHDp8OjbywgYXBsaWNhw6fDtWVzLjwvdzp0Pjx3OmJyLz48dzp0PjQpIEVuZXJnaWE6IHRyYWJhbGhvLCBmb3LDp2FzIGNvbnNlcnZhdGl2YXMsIGNvbnNlcnZhw6fDo28gZGUgZW5lcmdpYSBtZWPDom5pY2EsIGF0cml0bywgYXBsaWNhw6fDtWVzLjwvdzp0Pjx3OmJyLz48dzp0PjUpICBNb21lbnRvIGxpbmVhcjogY2VudHJvIGRlIG1hc3NhLCBzaXN0ZW1hIGRlIHBhcnTDrWN1bGFzLCBjb25zZXJ2YcOnw6NvIGRvIG1vbWVudG8gbGluZWFyLCBjb2xpc8O1ZXMsIGltcHVsc28uPC93OnQ+PHc6YnIvPjx3OnQ+NikgUm90YcOnw6NvOiB2YXJpw6F2ZWlzIGRvIG1vdmltZW50byByb3RhY2lvbmFsLCBlbmVyZ2lhIGNpbsOpdGljYSByb3RhY2lvbmFsLCBtb21lbnRvIGRlIGluw6lyY2lhLCB0b3JxdWUsIHJvbGFtZW50bywgY29uc2VydmHDp8OjbyBkbyBtb21lbnRvIGFuZ3VsYXIuPC93OnQ+PC93OnI+PC93OnA+"

# --- EN Programa (italic) ---
Replace-ParagraphWithXml "MSkgSW50cm9kdWN0aW9uIHRvIFBoeXNpY3M6IHVuaXQgc3lzdGVtcywgcmV2aWV3IG9mIHZlY3RvcnMsIGRpbWVuc2lvbmFsIGFuYWx5c2lzLiAyKSBLaW5lbWF0aWNzOiBvbmUgZGltZW5zaW9uYWwgbW90aW9uLCBmcmVlIGZhbGwsIGJpZGltZW5zaW9uYWwgbW90aW9uLCBwcm9qZWN0aWxlLiAgMykgRHluYW1pY3M6IE5ld3RvbuKAmXMgbGF3cywgZnJpY3Rpb24gZm9yY2UsIGRyYWcgZm9yY2UsIHRlcm1pbmFsIHNwZWVkLCB1bmlmb3JtIGNpcmN1bGFyIG1vdGlvbiwgZ3Jhdml0YXRpb24sIGFwcGxpY2F0aW9ucy40KSBFbmVyZ3k6IHdvcmssIGNvbnNlcnZhdGl2ZSBmb3JjZXMsIG1lY2hhbmljYWwgZW5lcmd5IGNvbnNlcnZhdGlvbiwgZnJpY3Rpb24sIGFwcGxpY2F0aW9ucy41KSAgTGluZWFyIG1vbWVudHVtOiBjZW50ZXIgb2YgbWFzcywgc3lzdGVtIG9mIHBhcnRpY2xlcywgY29uc2VydmF0aW9uIG9mIGxpbmVhciBtb21lbnR1bSwgY29sbGlzaW9ucywgaW1wdWxzZS42KSBSb3RhdGlvbjogcm90YXRpb25hbCB2YXJpYWJsZXMsIGtpbmV0aWMgZW5lcmd5IG9mIHJvdGF0aW9uLCByb3RhdGlvbmFsIGluZXJ0aWEsIHRvcnF1ZSwgcm9sbGluZywgY29uc2VydmF0aW9uIG9mIGFuZ3VsYXIgbW9tZW50dW0=" "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6clByPjx3OmkvPjwvdzpyUHI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4xKSBJbnRyb2R1Y3Rpb24gdG8gUGh5c2ljczogdW5pdCBzeXN0ZW1zLCByZXZpZXcgb2YgdmVjdG9ycywgZGltZW5zaW9uYWwgYW5hbHlzaXMuIDwvdzp0Pjx3OmJyLz48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPjIpIEtpbmVtYXRpY3M6IG9uZSBkaW1lbnNpb25hbCBtb3Rpb24sIGZyZWUgZmFsbCwgYmlkaW1lbnNpb25hbCBtb3Rpb24sIHByb2plY3RpbGUuIDwvdzp0Pjx3OmJyLz48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiAzKSBEeW5hbWljczogTmV3dG9u4oCZcyBsYXdzLCBmcmljdGlvbiBmb3JjZSwgZHJhZyBmb3JjZSwgdGVybWluYWwgc3BlZWQsIHVuaWZvcm0gY2lyY3VsYXIgbW90aW9uLCBncmF2aXRhdGlvbiwgYXBwbGljYXRpb25zLjwvdzp0Pjx3OmJyLz48dzp0PjQpIEVuZXJneTogd29yaywgY29uc2VydmF0aXZlIGZvcmNlcywgbWVjaGFuaWNhbCBlbmVyZ3kgY29uc2VydmF0aW9uLCBmcmljdGlvbiwgYXBwbGljYXRpb25zLjwvdzp0Pjx3OmJyLz48dzp0PjUpICBMaW5lYXIgbW9tZW50dW06IGNlbnRlciBvZiBtYXNzLCBzeXN0ZW0gb2YgcGFydGljbGVzLCBjb25zZXJ2YXRpb24gb2YgbGluZWFyIG1vbWVudHVtLCBjb2xsaXNpb25zLCBpbXB1bHNlLjwvdzp0Pjx3OmJyLz48dzp0PjYpIFJvdGF0aW9uOiByb3RhdGlvbmFsIHZhcmlhYmxlcywga2luZXRpYyBlbmVyZ3kgb2Ygcm90YXRpb24sIHJvdGF0aW9uYWwgaW5lcnRpYSwgdG9ycXVlLCByb2xsaW5nLCBjb25zZXJ2YXRpb24gb2YgYW5ndWxhciBtb21lbnR1bTwvdzp0PjwvdzpyPjwvdzpwPg=="

# --- Bibliografia ---
Replace-ParagraphWithXml "SEFMTElEQVksIEQ7IFJFU05JQ0ssIFIuIEZ1bmRhbWVudG9zIGRlIEbDrXNpY2EuIFZvbC4xLCBMVEMgKDIwMDgpLlNFQVJTLCBGLiBXLjsgWkVNQU5TS1ksIE0uIFcuOyBZT1VORywgSC4gRC47IEZSRUVETUFOLCBSLiBBLiBGw61zaWNhIEksIFZvbC4gMSwgUGVhcnNvbiBBZGRpc29uIFdlc2xleSAoMjAwOSkuSkVXRVRUIEpyLCBKb2huIFcuOyBTRVJXQVksIFJheW1vbmQgQS4gUHJpbmPDrXBpb3MgZGUgRsOtc2ljYS4gVm9sLiAxLCBUaG9tc29uIFBpb25laXJhICgyMDA4KS5OVVNTRU5aVkVJRywgSC5NLiBDdXJzbyBkZSBGw61zaWNhIELDoXNpY2EuIFZvbC4gMSwgRWRnYXJkIEJsdWNoZXIgKDIwMDgpLlRJUExFUiwgUC47IE1PU0NBLCBHLiBGw61zaWNhIHBhcmEgQ2llbnRpc3RhcyBlIEVuZ2VuaGVpcm9zLiBWb2wuMSwgTFRDICgyMDA4KS4=" "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6dD5IQUxMSURBWSwgRDsgUkVTTklDSywgUi4gRnVuZGFtZW50b3MgZGUgRsOtc2ljYS4gVm9sLjEsIExUQyAoMjAwOCkuPC93OnQ+PHc6YnIvPjx3OnQ+U0VBUlMsIEYuIFcuOyBaRU1BTlNLWSwgTS4gVy47IFlPVU5HLCBILiBELjsgRlJFRURNQU4sIFIuIEEuIEbDrXNpY2EgSSwgVm9sLiAxLCBQZWFyc29uIEFkZGlzb24gV2VzbGV5ICgyMDA5KS48L3c6dD48dzpici8+PHc6dD5KRVdFVFQgSnIsIEpvaG4gVy47IFNFUldBWSwgUmF5bW9uZCBBLiBQcmluY8OtcGlvcyBkZSBGw61zaWNhLiBWb2wuIDEsIFRob21zb24gUGlvbmVpcmEgKDIwMDgpLjwvdzp0Pjx3OmJyLz48dzp0Pk5VU1NFTlpWRUlHLCBILk0uIEN1cnNvIGRlIEbDrXNpY2EgQsOhc2ljYS4gVm9sLiAxLCBFZGdhcmQgQmx1Y2hlciAoMjAwOCkuPC93OnQ+PHc6YnIvPjx3OnQ+VElQTEVSLCBQLjsgTU9TQ0EsIEcuIEbDrXNpY2EgcGFyYSBDaWVudGlzdGFzIGUgRW5nZW5oZWlyb3MuIFZvbC4xLCBMVEMgKDIwMDgpLjwvdzp0PjwvdzpyPjwvdzpwPg=="

Write-Host "Done splitting paragraphs."
